$d = $word.ActiveDocument

# Locate the unique heading paragraph "Estensione 1c" that immediately
# follows the table we need to insert a new (empty) paragraph before.
$headingRange = $d.Content
$headingRange.Find.Execute("Estensione 1c", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headingPara = $headingRange.Paragraphs(1)
$headingStart = $headingPara.Range.Start

# Collapsed insertion point right before the heading paragraph (i.e. right
# after the preceding table's closing </w:tbl>).
$insertionPoint = $d.Range($headingStart, $headingStart)

# Insert a new, empty paragraph with the same "spacer heading" styling used
# elsewhere in the document, but with a near-invisible (2 half-points) end
# of paragraph mark font size, matching the target markup exactly.
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="400" w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:outlineLvl w:val="0"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="00615E"/><w:kern w:val="36"/><w:sz w:val="2"/><w:szCs w:val="2"/><w:lang w:eastAsia="it-IT"/></w:rPr></w:pPr></w:p>'

$insertionPoint.InsertXML($newParaXml) | Out-Null

Write-Output "Inserted spacer paragraph before 'Estensione 1c'."
